$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" under duplicate_image_filename (column E) for rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
